# edit.ps1 - Applies the "Updated slide, added scope, changed title" commit
# to the Project 1 PPT presentation.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 (Title slide): fix the TITLE line text.
# ---------------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
$subtitle = $titleSlide.Shapes.Item(2)
$titlePara = $subtitle.TextFrame.TextRange.Paragraphs(2, 1)
$titlePara.Text = "TITLE: The Impact of Covid-19 on Suicide Rate"

# ---------------------------------------------------------------------------
# 2) Slide 2 (Introduction): remove the second paragraph ("The goal of this
#    project is to examine the impact of Covid-19 on suicide rate.") from the
#    content placeholder - that idea is moved to the new Scope slide.
# ---------------------------------------------------------------------------
$introSlide = $p.Slides.Item(2)
$introBody = $introSlide.Shapes.Item(2)
$introBody.TextFrame.TextRange.Text = "Covid-19, a world pandemic which broke out in late 2019 brought significant changes to the world. Despite its direct impact on sickness and mental health around the world, its effect was also felt on suicide rate. "

# ---------------------------------------------------------------------------
# 3) Insert a brand-new "Scope" slide as slide 3 (Title and Content layout),
#    pushing Data Source/Observation/Code/Code continued down by one.
# ---------------------------------------------------------------------------
$contentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$scopeSlide = $p.Slides.AddSlide(3, $contentLayout)

$scopeTitle = $scopeSlide.Shapes.Item(1)
$scopeTitle.TextFrame.TextRange.Text = "Scope"

$scopeBody = $scopeSlide.Shapes.Item(2)
$scopeBodyTf = $scopeBody.TextFrame
$scopeBodyTf.TextRange.Text = "The goal of this project is:`r To examine the impact of Covid-19 on suicide rate. The task includes gathering information on global suicide rate and recorded Covid-19 cases.`rThe information will be broken down using different strategies to identify and point out any patterns.`r"
